$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.223867666666667
$ws.Range("H2").Value = 3.671603
$ws.Range("I2").Value = 0.2362882870487179
$ws.Range("J2").Value = 0.2492245847267186
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.67754233333333
$ws.Range("N2").Value = 50.03262699999999
$ws.Range("O2").Value = 0.9535192900707901
$ws.Range("P2").Value = 0.9578676752791928
$ws.Range("Q2").Value = 20.41110482123122
$ws.Range("R2").Value = 183.699943391081
$ws.Range("S2").Value = 0.2253054397187365
$ws.Range("T2").Value = 0.2387241735946042

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.223867666666667
$ws.Range("H3").Value = 3.671603
$ws.Range("I3").Value = 0.2362882870487179
$ws.Range("J3").Value = 0.2492245847267186
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.574769
$ws.Range("N3").Value = 1.724307
$ws.Range("O3").Value = 0.0328617561197435
$ws.Range("P3").Value = 0.03301161735036698
$ws.Range("Q3").Value = 0.7034411949023334
$ws.Range("R3").Value = 6.330970754121
$ws.Range("S3").Value = 0.007764848062946913
$ws.Range("T3").Value = 0.00822730662530255

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.223867666666667
$ws.Range("H4").Value = 3.671603
$ws.Range("I4").Value = 0.2362882870487179
$ws.Range("J4").Value = 0.2492245847267186
$ws.Range("K4").Value = 2
$ws.Range("M4").Value = 0.2382025
$ws.Range("N4").Value = 0.476405
$ws.Range("O4").Value = 0.01361895380946642
$ws.Range("P4").Value = 0.009120707370440172
$ws.Range("Q4").Value = 0.2915283378691667
$ws.Range("R4").Value = 1.749170027215
$ws.Range("S4").Value = 0.00321799926703443
$ws.Range("T4").Value = 0.002273104506811874

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.9252306666666666
$ws.Range("H5").Value = 2.775692
$ws.Range("I5").Value = 0.178631379278977
$ws.Range("J5").Value = 0.1884110798551137
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.67754233333333
$ws.Range("N5").Value = 50.03262699999999
$ws.Range("O5").Value = 0.9535192900707901
$ws.Range("P5").Value = 0.9578676752791928
$ws.Range("Q5").Value = 15.43057361143155
$ws.Range("R5").Value = 138.875162502884
$ws.Range("S5").Value = 0.1703284659544562
$ws.Range("T5").Value = 0.1804728830576601

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9252306666666666
$ws.Range("H6").Value = 2.775692
$ws.Range("I6").Value = 0.178631379278977
$ws.Range("J6").Value = 0.1884110798551137
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.574769
$ws.Range("N6").Value = 1.724307
$ws.Range("O6").Value = 0.0328617561197435
$ws.Range("P6").Value = 0.03301161735036698
$ws.Range("Q6").Value = 0.5317939050493333
$ws.Range("R6").Value = 4.786145145443999
$ws.Range("S6").Value = 0.005870140821199144
$ws.Range("T6").Value = 0.006219754472746449

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.9252306666666666
$ws.Range("H7").Value = 2.775692
$ws.Range("I7").Value = 0.178631379278977
$ws.Range("J7").Value = 0.1884110798551137
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.2382025
$ws.Range("N7").Value = 0.476405
$ws.Range("O7").Value = 0.01361895380946642
$ws.Range("P7").Value = 0.009120707370440172
$ws.Range("Q7").Value = 0.2203922578766667
$ws.Range("R7").Value = 1.32235354726
$ws.Range("S7").Value = 0.002432772503321664
$ws.Range("T7").Value = 0.001718442324707127

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.8203320000000001
$ws.Range("H8").Value = 2.460996
$ws.Range("I8").Value = 0.1583789231226106
$ws.Range("J8").Value = 0.1670498433828809
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.67754233333333
$ws.Range("N8").Value = 50.03262699999999
$ws.Range("O8").Value = 0.9535192900707901
$ws.Range("P8").Value = 0.9578676752791928
$ws.Range("Q8").Value = 13.681121657388
$ws.Range("R8").Value = 123.130094916492
$ws.Range("S8").Value = 0.1510173583380479
$ws.Range("T8").Value = 0.1600116451369134

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.8203320000000001
$ws.Range("H9").Value = 2.460996
$ws.Range("I9").Value = 0.1583789231226106
$ws.Range("J9").Value = 0.1670498433828809
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.574769
$ws.Range("N9").Value = 1.724307
$ws.Range("O9").Value = 0.0328617561197435
$ws.Range("P9").Value = 0.03301161735036698
$ws.Range("Q9").Value = 0.471501403308
$ws.Range("R9").Value = 4.243512629772001
$ws.Range("S9").Value = 0.005204609546162834
$ws.Range("T9").Value = 0.005514585508194398

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8203320000000001
$ws.Range("H10").Value = 2.460996
$ws.Range("I10").Value = 0.1583789231226106
$ws.Range("J10").Value = 0.1670498433828809
$ws.Range("K10").Value = 2
$ws.Range("M10").Value = 0.2382025
$ws.Range("N10").Value = 0.476405
$ws.Range("O10").Value = 0.01361895380946642
$ws.Range("P10").Value = 0.009120707370440172
$ws.Range("Q10").Value = 0.19540513323
$ws.Range("R10").Value = 1.17243079938
$ws.Range("S10").Value = 0.002156955238399866
$ws.Range("T10").Value = 0.001523612737773118

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.40357
$ws.Range("H11").Value = 4.210710000000001
$ws.Range("I11").Value = 0.2709828522198361
$ws.Range("J11").Value = 0.285818605975276
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 16.67754233333333
$ws.Range("N11").Value = 50.03262699999999
$ws.Range("O11").Value = 0.9535192900707901
$ws.Range("P11").Value = 0.9578676752791928
$ws.Range("Q11").Value = 23.40809809279667
$ws.Range("R11").Value = 210.67288283517
$ws.Range("S11").Value = 0.258387376870016
$ws.Range("T11").Value = 0.2737764036570773

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.40357
$ws.Range("H12").Value = 4.210710000000001
$ws.Range("I12").Value = 0.2709828522198361
$ws.Range("J12").Value = 0.285818605975276
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.574769
$ws.Range("N12").Value = 1.724307
$ws.Range("O12").Value = 0.0328617561197435
$ws.Range("P12").Value = 0.03301161735036698
$ws.Range("Q12").Value = 0.8067285253300001
$ws.Range("R12").Value = 7.260556727970001
$ws.Range("S12").Value = 0.008904972402280747
$ws.Range("T12").Value = 0.009435334452071125

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.40357
$ws.Range("H13").Value = 4.210710000000001
$ws.Range("I13").Value = 0.2709828522198361
$ws.Range("J13").Value = 0.285818605975276
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.2382025
$ws.Range("N13").Value = 0.476405
$ws.Range("O13").Value = 0.01361895380946642
$ws.Range("P13").Value = 0.009120707370440172
$ws.Range("Q13").Value = 0.334333882925
$ws.Range("R13").Value = 2.00600329755
$ws.Range("S13").Value = 0.003690502947539412
$ws.Range("T13").Value = 0.002606867866127635

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.8065525
$ws.Range("H14").Value = 1.613105
$ws.Range("I14").Value = 0.1557185583298584
$ws.Range("J14").Value = 0.1094958860600107
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 16.67754233333333
$ws.Range("N14").Value = 50.03262699999999
$ws.Range("O14").Value = 0.9535192900707901
$ws.Range("P14").Value = 0.9578676752791928
$ws.Range("Q14").Value = 13.45131346280583
$ws.Range("R14").Value = 80.70788077683498
$ws.Range("S14").Value = 0.1484806491895335
$ws.Range("T14").Value = 0.1048825698329378

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.8065525
$ws.Range("H15").Value = 1.613105
$ws.Range("I15").Value = 0.1557185583298584
$ws.Range("J15").Value = 0.1094958860600107
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.574769
$ws.Range("N15").Value = 1.724307
$ws.Range("O15").Value = 0.0328617561197435
$ws.Range("P15").Value = 0.03301161735036698
$ws.Range("Q15").Value = 0.4635813738725
$ws.Range("R15").Value = 2.781488243235
$ws.Range("S15").Value = 0.005117185287153859
$ws.Range("T15").Value = 0.003614636292052455

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.8065525
$ws.Range("H16").Value = 1.613105
$ws.Range("I16").Value = 0.1557185583298584
$ws.Range("J16").Value = 0.1094958860600107
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 0.2382025
$ws.Range("N16").Value = 0.476405
$ws.Range("O16").Value = 0.01361895380946642
$ws.Range("P16").Value = 0.009120707370440172
$ws.Range("Q16").Value = 0.19212282188125
$ws.Range("R16").Value = 0.7684912875250001
$ws.Range("S16").Value = 0.002120723853171043
$ws.Range("T16").Value = 0.0009986799350204167
